$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(8).Delete()
$ws.Range("A2:J2").Select()
